# Fix trend line baseline to use actual mean at center year
# Collapses the table from (Variable, Class, Slope, Intercept, R_squared,
# p_value, n_years) down to (Variable, Slope, p_value) and refreshes the
# Slope / p_value numbers (and reorders two rows) to reflect the corrected
# regression baseline.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row -------------------------------------------------------
$ws.Range("B1").Value = "Slope"
$ws.Range("C1").Value = "p_value"

# --- Row 2: Veteran Status -------------------------------------------------
$ws.Range("A2").Value = "Veteran Status"
$ws.Range("B2").Value = -0.134841301067288
$ws.Range("C2").Value = 0.00135846697049377

# --- Row 3: Family Income ---------------------------------------------------
$ws.Range("A3").Value = "Family Income"
$ws.Range("B3").Value = 0.2584828359824106
$ws.Range("C3").Value = 0.0192849003679574

# --- Row 4: Voting Method (was State Representative) ------------------------
$ws.Range("A4").Value = "Voting Method"
$ws.Range("B4").Value = 0.137223807092183
$ws.Range("C4").Value = 0.04497397374979856

# --- Row 5: State Representative (was Voting Method) -------------------------
$ws.Range("A5").Value = "State Representative"
$ws.Range("B5").Value = 0.1116634092427386
$ws.Range("C5").Value = 0.1229602268622791

# --- Row 6: Employment Status -----------------------------------------------
$ws.Range("A6").Value = "Employment Status"
$ws.Range("B6").Value = 0.2152019105524243
$ws.Range("C6").Value = 0.1486843722550557

# --- Row 7: Union Membership -------------------------------------------------
$ws.Range("A7").Value = "Union Membership"
$ws.Range("B7").Value = 0.08367787393634148
$ws.Range("C7").Value = 0.1494887215625006

# --- Row 8: U.S. House --------------------------------------------------------
$ws.Range("A8").Value = "U.S. House"
$ws.Range("B8").Value = 0.09288080445776516
$ws.Range("C8").Value = 0.3989869415835802

# --- Row 9: State Senator ------------------------------------------------------
$ws.Range("A9").Value = "State Senator"
$ws.Range("B9").Value = -0.02002032986058871
$ws.Range("C9").Value = 0.8541625635049479

# --- Drop the now-unused Intercept / R_squared / p_value(old col F) / n_years
#     columns entirely so the used range shrinks back to A1:C9. ----------------
$ws.Range("D1:G9").Clear()
